$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.956.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.884.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7432'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +0.84%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07242'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.62%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.135.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.73%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08349'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7568'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.410'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.152'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.074.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.97%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '248.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007876'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.163.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9979'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.049'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9987'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1560'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.320'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.045'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.507'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.607'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.538'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.230'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05370'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.53%  '

$ws.Range("E35").Value = '  +0.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7572'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.007'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.708'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01969'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.761'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4565'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.00%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.112.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.88%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.158'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8618'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '105.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.14%  '

$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.877'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.628'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.073.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.81%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.568'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.88%  '
